$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 48: TRADING_ATTEMPT for SOL
$ws.Cells.Item(48, 1).Value = "2025-09-25T12:43:56.425956"
$ws.Cells.Item(48, 2).Value = "TRADING_ATTEMPT"
$ws.Cells.Item(48, 3).Value = "SOL"
$ws.Cells.Item(48, 4).Value = "UNKNOWN"
$ws.Cells.Item(48, 5).Value = 211.5960881268497
$ws.Cells.Item(48, 11).Value = "ATTEMPT"
$ws.Cells.Item(48, 12).Value = "Attempting trade 1/1"

# Row 49: POSITION_FAILED for SOL
$ws.Cells.Item(49, 1).Value = "2025-09-25T12:43:56.817993"
$ws.Cells.Item(49, 2).Value = "POSITION_FAILED"
$ws.Cells.Item(49, 3).Value = "SOL"
$ws.Cells.Item(49, 4).Value = "UNKNOWN"
$ws.Cells.Item(49, 11).Value = "FAILED"
$ws.Cells.Item(49, 12).Value = "Trade execution failed for trade 1"
